$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H9").Value = 3921.3333
$ws.Range("J9").Value = 17409.5
$ws.Range("L9").Value = 17409.5
$ws.Range("N9").Value = -17747.5
$ws.Range("H18").Value = 1366.2
$ws.Range("I18").Value = 963.8570999999999
$ws.Range("J18").Value = 6999
$ws.Range("K18").Value = 963.8570999999999
$ws.Range("L18").Value = 6999
$ws.Range("M18").Value = -679.8570999999999
$ws.Range("N18").Value = -7567
$ws.Range("H19").Value = 1722.625
$ws.Range("I19").Value = 2036.3846
$ws.Range("J19").Value = 1351.8182
$ws.Range("K19").Value = 2036.3846
$ws.Range("L19").Value = 1351.8182
$ws.Range("M19").Value = -1861.3846
$ws.Range("N19").Value = -1701.8182
$ws.Range("H38").Value = 652
$ws.Range("I38").Value = 182.4
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 547.2
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -175.2
$ws.Range("N38").Value = -9744
$ws.Range("H106").Value = 4032.2666
$ws.Range("I106").Value = 3407.9092
$ws.Range("K106").Value = 3407.9092
$ws.Range("M106").Value = -2776.9092
$ws.Range("H137").Value = 2206.5715
$ws.Range("I137").Value = 1078.4166
$ws.Range("K137").Value = 3235.2498
$ws.Range("M137").Value = -685.2498000000001
$ws.Range("H138").Value = 5613.24
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 5613.24
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 16839.72
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -27119.72

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H32").Value = 1266.7567
$ws.Range("I32").Value = 1090.0303
$ws.Range("K32").Value = 1090.0303
$ws.Range("M32").Value = -803.0302999999999
$ws.Range("H74").Value = 2392.4324
$ws.Range("I74").Value = 1080.871
$ws.Range("K74").Value = 1080.871
$ws.Range("M74").Value = -206.8710000000001
$ws.Range("H77").Value = 2392.4324
$ws.Range("I77").Value = 1080.871
$ws.Range("K77").Value = 5404.355
$ws.Range("M77").Value = -1036.355
$ws.Range("H97").Value = 991.9048
$ws.Range("I97").Value = 924.2353000000001
$ws.Range("K97").Value = 924.2353000000001
$ws.Range("M97").Value = -428.2353000000001
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H29").Value = 349
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H94").Value = 619.36365
$ws.Range("I94").Value = 457.55554
$ws.Range("K94").Value = 457.55554
$ws.Range("M94").Value = -6.555540000000008
$ws.Range("H96").Value = 12400.5
$ws.Range("I96").Value = 12400.5
$ws.Range("K96").Value = 12400.5
$ws.Range("M96").Value = -9654.5
$ws.Range("H99").Value = 29232.842
$ws.Range("I99").Value = 34464.062
$ws.Range("K99").Value = 34464.062
$ws.Range("M99").Value = -32966.062

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H12").Value = 1565.625
$ws.Range("I12").Value = 105.2
$ws.Range("J12").Value = 3999.6667
$ws.Range("K12").Value = 105.2
$ws.Range("L12").Value = 3999.6667
$ws.Range("M12").Value = 64.8
$ws.Range("N12").Value = -4339.6667
$ws.Range("H31").Value = 13377.737
$ws.Range("I31").Value = 1114.3182
$ws.Range("K31").Value = 1114.3182
$ws.Range("M31").Value = -819.3181999999999
$ws.Range("H34").Value = 13377.737
$ws.Range("I34").Value = 1114.3182
$ws.Range("K34").Value = 1114.3182
$ws.Range("M34").Value = -912.3181999999999
$ws.Range("H69").Value = 55750
$ws.Range("I69").Value = 97500
$ws.Range("J69").Value = 14000
$ws.Range("K69").Value = 97500
$ws.Range("L69").Value = 14000
$ws.Range("M69").Value = -96751
$ws.Range("N69").Value = -15498
$ws.Range("H72").Value = 55750
$ws.Range("I72").Value = 97500
$ws.Range("J72").Value = 14000
$ws.Range("K72").Value = 292500
$ws.Range("L72").Value = 42000
$ws.Range("M72").Value = -288756
$ws.Range("N72").Value = -49488
$ws.Range("H99").Value = 3858.4
$ws.Range("I99").Value = 4323
$ws.Range("K99").Value = 4323
$ws.Range("M99").Value = -2825
$ws.Range("H126").Value = 3858.4
$ws.Range("I126").Value = 4323
$ws.Range("K126").Value = 12969
$ws.Range("M126").Value = -10499

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H40").Value = 116.5
$ws.Range("I40").Value = 99
$ws.Range("J40").Value = 120
$ws.Range("K40").Value = 396
$ws.Range("L40").Value = 480
$ws.Range("M40").Value = -327
$ws.Range("N40").Value = -618
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H107").Value = 689.1429000000001
$ws.Range("I107").Value = 397.66666
$ws.Range("J107").Value = 907.75
$ws.Range("K107").Value = 1192.99998
$ws.Range("L107").Value = 2723.25
$ws.Range("M107").Value = 727.0000199999999
$ws.Range("N107").Value = -6563.25
$ws.Range("H132").Value = 1124.25
$ws.Range("J132").Value = 1498.5
$ws.Range("L132").Value = 13486.5
$ws.Range("N132").Value = -18546.5
$ws.Range("H136").Value = 5064.5557
$ws.Range("I136").Value = 4086.1667
$ws.Range("K136").Value = 12258.5001
$ws.Range("M136").Value = -7158.500100000001
$ws.Range("H137").Value = 5378.375
$ws.Range("I137").Value = 1259.5
$ws.Range("J137").Value = 9497.25
$ws.Range("K137").Value = 3778.5
$ws.Range("L137").Value = 28491.75
$ws.Range("M137").Value = 1321.5
$ws.Range("N137").Value = -38691.75
$ws.Range("H140").Value = 3234.4211
$ws.Range("I140").Value = 1438
$ws.Range("J140").Value = 4063.5386
$ws.Range("K140").Value = 4314
$ws.Range("L140").Value = 12190.6158
$ws.Range("M140").Value = 866
$ws.Range("N140").Value = -22550.6158

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H18").Value = 8000
$ws.Range("J18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("N18").Value = -8586
$ws.Range("H53").Value = 5620
$ws.Range("I53").Value = 5200
$ws.Range("J53").Value = 5725
$ws.Range("K53").Value = 5200
$ws.Range("L53").Value = 5725
$ws.Range("M53").Value = -4569
$ws.Range("N53").Value = -6987
$ws.Range("H102").Value = 2751.9312
$ws.Range("I102").Value = 2012.2858
$ws.Range("J102").Value = 3442.2666
$ws.Range("K102").Value = 2012.2858
$ws.Range("L102").Value = 3442.2666
$ws.Range("M102").Value = -390.2858000000001
$ws.Range("N102").Value = -6686.2666
$ws.Range("H132").Value = 9461.333000000001
$ws.Range("I132").Value = 10630.5
$ws.Range("J132").Value = 7999.875
$ws.Range("K132").Value = 31891.5
$ws.Range("L132").Value = 23999.625
$ws.Range("M132").Value = -29361.5
$ws.Range("N132").Value = -29059.625

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 3033.6667
$ws.Range("I46").Value = 1431.3334
$ws.Range("J46").Value = 3674.6
$ws.Range("K46").Value = 1431.3334
$ws.Range("L46").Value = 3674.6
$ws.Range("M46").Value = -1243.3334
$ws.Range("N46").Value = -4050.6
$ws.Range("H61").Value = 3800
$ws.Range("I61").Value = 3560
$ws.Range("K61").Value = 3560
$ws.Range("M61").Value = -3358
$ws.Range("H100").Value = 9318.429
$ws.Range("I100").Value = 2051.4546
$ws.Range("J100").Value = 35964
$ws.Range("K100").Value = 2051.4546
$ws.Range("L100").Value = 35964
$ws.Range("M100").Value = -1510.4546
$ws.Range("N100").Value = -37046
$ws.Range("H113").Value = 3800
$ws.Range("I113").Value = 3560
$ws.Range("K113").Value = 3560
$ws.Range("M113").Value = -1390

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H99").Value = 75001
$ws.Range("I99").Value = 75001
$ws.Range("K99").Value = 75001
$ws.Range("M99").Value = -72006
$ws.Range("H122").Value = 1497.6451
$ws.Range("I122").Value = 1574.4073
$ws.Range("J122").Value = 979.5
$ws.Range("K122").Value = 4723.2219
$ws.Range("L122").Value = 2938.5
$ws.Range("M122").Value = -2273.2219
$ws.Range("N122").Value = -7838.5
$ws.Range("H126").Value = 3147.423
$ws.Range("J126").Value = 4263.222
$ws.Range("L126").Value = 12789.666
$ws.Range("N126").Value = -17729.666
$ws.Range("H132").Value = 5511.3335
$ws.Range("I132").Value = 1433.5
$ws.Range("J132").Value = 13667
$ws.Range("K132").Value = 4300.5
$ws.Range("L132").Value = 41001
$ws.Range("M132").Value = -1770.5
$ws.Range("N132").Value = -46061
$ws.Range("H136").Value = 6701.3335
$ws.Range("I136").Value = 2071.625
$ws.Range("K136").Value = 6214.875
$ws.Range("M136").Value = -3664.875
